{"js": "// Update the date line and all \"N\u00f7N=\" division problems in the table.\nconst replacements = [\n  [\"2024-06-07 Friday\", \"2024-06-08 Saturday\"],\n  [\"856\u00f78=\", \"612\u00f76=\"],\n  [\"329\u00f79=\", \"921\u00f76=\"],\n  [\"834\u00f78=\", \"349\u00f75=\"],\n  [\"880\u00f74=\", \"384\u00f79=\"],\n  [\"933\u00f77=\", \"850\u00f78=\"],\n  [\"329\u00f77=\", \"341\u00f78=\"],\n  [\"602\u00f73=\", \"530\u00f78=\"],\n  [\"986\u00f73=\", \"756\u00f77=\"],\n  [\"646\u00f78=\", \"491\u00f76=\"],\n  [\"611\u00f78=\", \"537\u00f76=\"],\n  [\"491\u00f77=\", \"536\u00f73=\"],\n  [\"285\u00f73=\", \"692\u00f73=\"],\n  [\"177\u00f73=\", \"573\u00f78=\"],\n  [\"558\u00f74=\", \"583\u00f72=\"],\n  [\"843\u00f79=\", \"684\u00f73=\"],\n  [\"130\u00f72=\", \"664\u00f72=\"],\n  [\"521\u00f74=\", \"290\u00f74=\"],\n  [\"379\u00f73=\", \"695\u00f78=\"],\n  [\"418\u00f74=\", \"767\u00f76=\"],\n  [\"383\u00f72=\", \"170\u00f75=\"],\n  [\"899\u00f77=\", \"798\u00f77=\"],\n  [\"792\u00f78=\", \"304\u00f75=\"],\n  [\"696\u00f77=\", \"478\u00f74=\"],\n  [\"802\u00f77=\", \"847\u00f78=\"],\n  [\"958\u00f76=\", \"443\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2024-06-07 Friday\", \"2024-06-08 Saturday\"),\n  @(\"856\u00f78=\", \"612\u00f76=\"),\n  @(\"329\u00f79=\", \"921\u00f76=\"),\n  @(\"834\u00f78=\", \"349\u00f75=\"),\n  @(\"880\u00f74=\", \"384\u00f79=\"),\n  @(\"933\u00f77=\", \"850\u00f78=\"),\n  @(\"329\u00f77=\", \"341\u00f78=\"),\n  @(\"602\u00f73=\", \"530\u00f78=\"),\n  @(\"986\u00f73=\", \"756\u00f77=\"),\n  @(\"646\u00f78=\", \"491\u00f76=\"),\n  @(\"611\u00f78=\", \"537\u00f76=\"),\n  @(\"491\u00f77=\", \"536\u00f73=\"),\n  @(\"285\u00f73=\", \"692\u00f73=\"),\n  @(\"177\u00f73=\", \"573\u00f78=\"),\n  @(\"558\u00f74=\", \"583\u00f72=\"),\n  @(\"843\u00f79=\", \"684\u00f73=\"),\n  @(\"130\u00f72=\", \"664\u00f72=\"),\n  @(\"521\u00f74=\", \"290\u00f74=\"),\n  @(\"379\u00f73=\", \"695\u00f78=\"),\n  @(\"418\u00f74=\", \"767\u00f76=\"),\n  @(\"383\u00f72=\", \"170\u00f75=\"),\n  @(\"899\u00f77=\", \"798\u00f77=\"),\n  @(\"792\u00f78=\", \"304\u00f75=\"),\n  @(\"696\u00f77=\", \"478\u00f74=\"),\n  @(\"802\u00f77=\", \"847\u00f78=\"),\n  @(\"958\u00f76=\", \"443\u00f74=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
